# edit.ps1 -- applies the yearly_report.docx revision described by the
# commit "feat: final version for deployment".
#
# Strategy:
#   1. Text-only paragraph rewrites are done with Find/Replace (content is
#      unique, so this is immune to paragraph-index shifting).
#   2. Structural changes (paragraphs added/removed) are done through
#      $d.Paragraphs(i).Range, walking the document from the END towards
#      the START (by the paragraph's ORIGINAL 1-based index) so that an
#      insertion/deletion never invalidates an index used by a later
#      (i.e. earlier-in-the-document) step.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# ---------------------------------------------------------------------
# 1. Plain text replacements (no paragraph count change)
# ---------------------------------------------------------------------

Replace-Text "Dear Friends, Partners, and Supporters," "Dear Friends and Supporters of the IKEA Foundation,"

Replace-Text "As we reflect on the year 2023, I am filled with gratitude and pride for the remarkable strides we have made together. This year has been one of both significant achievements and formidable challenges, each shaping our journey towards a more sustainable and equitable world." "As we reflect on the year 2023, it is with a profound sense of gratitude and commitment that I address you in this annual report. This year has been a remarkable journey filled with substantial achievements and, admittedly, some challenges that have tested our resilience and strengthened our resolve."

Replace-Text "The IKEA Foundation remains steadfast in its commitment to improving the lives of the many people, particularly those facing poverty and the adverse impacts of climate change. Our initiatives across the globe, in partnership with over 140 dedicated organizations, have aimed to transform lives and communities sustainably." "The IKEA Foundation has continued to make significant strides in our mission to create a better everyday life for the many people, particularly those living in vulnerable conditions. Our efforts this year have been notably impactful in addressing the dual threats of poverty and climate change, which remain at the forefront of our work."

Replace-Text "Achievements and Challenges" "Achievements and Progress"

Replace-Text 'This year, we have seen impactful progress in various sectors, notably in our efforts to combat climate change and enhance economic growth within vulnerable communities. Our project, "The Funders Table," which aligns with the United Nations'' Goal 13, has made significant headway in reducing carbon emissions and fostering global cooperation against climate change[x]. Similarly, the "Just Transition Fund" has been pivotal in promoting the use of renewable energy and sustainable practices in South Africa, Vietnam, and Indonesia[x].' "This year, we saw the successful initiation and continuation of key projects such as the African Private Sector Forum on Forced Displacement [x], which began in May 2023. This project, based in Nairobi, Kenya, is a testament to our commitment to promoting Decent Work and Economic Growth as per the United Nations' Goal 8 [x]. With a total budget of `$500,000 [x], and a significant portion already funded, this project has been making strides in increasing employment opportunities by 30% within the refugee communities [x]."

Replace-Text 'However, the journey was not without its hurdles. Political instability and logistical challenges in regions like Nairobi and Northern Syria have tested our resilience and adaptability. The "Protracted Displacement in an Urban Context" project faced difficulties due to the complex dynamics of urban environments and political uncertainties[x]. Despite these challenges, our team''s dedication and our partners'' support have been instrumental in navigating these complexities and continuing our mission.' "In the realm of environmental sustainability, our project, The Funders Table, has been a catalyst in combating climate change. This initiative, which aligns with United Nations' Goal 13 [x], commenced in April 2023 with a robust budget of `$5,000,000 [x]. Our foundation's contribution of 40% towards this project underscores our commitment to reducing carbon emissions by 30% within the funded projects [x]."

Replace-Text "Gratitude and Future Outlook" "Challenges and Learning"

Replace-Text "I extend my deepest thanks to you, our partners and supporters, whose unwavering commitment has been the cornerstone of our achievements. Your collaboration, insights, and generous contributions have been crucial in bringing our shared vision to life." "While we celebrate our achievements, we also acknowledge the challenges faced along the way. Political instability and logistical issues in certain regions have posed risks to the timely progress of our projects. However, these challenges have provided us with valuable insights and learning opportunities to enhance our strategies and operations."

Replace-Text "Conclusion" "Looking Forward"

Replace-Text "In closing, let us celebrate the successes of the past year and brace for the opportunities and challenges ahead. Together, we are making a substantial difference in the lives of many and contributing to a healthier, more sustainable planet. Here's to continuing our journey with hope, determination, and shared resolve." "As we look to the future, our vision for the next year and beyond is to not only continue these impactful projects but also to innovate further in both poverty alleviation and climate action. We aim to expand our geographical reach and deepen our impact through strategic partnerships and enhanced project models."

Replace-Text "Thank you for being a pivotal part of our journey." "Thank you for your continued support and belief in our mission. Together, we are making a difference."

Replace-Text "Chairperson/President, IKEA Foundation" "Chairperson/President"

# ---------------------------------------------------------------------
# 2. Structural changes -- walk from the bottom of the document upward
#    so earlier (original) paragraph indices stay valid.
# ---------------------------------------------------------------------

# --- (a) Drop the "Contact Information" + "Note" block (orig paragraphs
#         32-36), then turn the new last paragraph into "IKEA Foundation [x]".
$d.Paragraphs(36).Range.Delete()   # "Note: ..."
$d.Paragraphs(35).Range.Delete()   # empty
$d.Paragraphs(34).Range.Delete()   # "For inquiries ..."
$d.Paragraphs(33).Range.Delete()   # "Contact Information:"
$d.Paragraphs(32).Range.Delete()   # empty

$d.Paragraphs(32).Range.Text = "IKEA Foundation [x]"

# --- (b) Remove the "Looking ahead ..." paragraph (orig 20) together with
#         its leading blank separator (orig 19).
$d.Paragraphs(20).Range.Delete()
$d.Paragraphs(19).Range.Delete()

# --- (c) After "In closing ..." (orig paragraph 24, already retexted to the
#         "As we look to the future ..." copy above), insert the new
#         "Gratitude" section: blank / heading / blank / paragraph.
$d.Paragraphs(24).Range.InsertParagraphAfter()
$d.Paragraphs(25).Range.InsertParagraphAfter()
$d.Paragraphs(26).Range.InsertParagraphAfter()
$d.Paragraphs(27).Range.InsertParagraphAfter()

$p26 = $d.Paragraphs(26)
$p26.Range.Text = "Gratitude"
$p26.Range.Font.Bold = 1

$d.Paragraphs(28).Range.Text = "I want to express my deepest gratitude to all our partners, donors, and team members whose relentless dedication and support have been indispensable. It is your commitment that empowers us to keep striving towards a better world."

# --- (d) After "However, the journey ..." (orig paragraph 14, already
#         retexted to the "In the realm of environmental sustainability ..."
#         copy above), insert the new "Just Transition Fund" paragraph.
$d.Paragraphs(14).Range.InsertParagraphAfter()
$d.Paragraphs(15).Range.InsertParagraphAfter()

$d.Paragraphs(16).Range.Text = "Another notable project is the Just Transition Fund, which began in April 2023. With a focus on Climate Action, this global initiative aims to reduce carbon emissions by 25% and increase renewable energy use by 40% in countries like South Africa, Vietnam, and Indonesia [x]. Although it faces challenges, the project's strong partnerships with influential global organizations have been pivotal in navigating the complexities of regulatory changes across multiple countries [x]."

Write-Output "done"
